$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 1
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 109

$ws.Range("H8").Value = 41.4
$ws.Range("I8").Value = 41.4
$ws.Range("K8").Value = 124.2
$ws.Range("M8").Value = 14.80000000000001

$ws.Range("H15").Value = 1444.8667
$ws.Range("I15").Value = 1444.8667
$ws.Range("K15").Value = 4334.6001
$ws.Range("M15").Value = -4165.6001

$ws.Range("H49").Value = 98.333336
$ws.Range("J49").Value = 48
$ws.Range("L49").Value = 144
$ws.Range("N49").Value = -416

$ws.Range("H51").Value = 29999
$ws.Range("J51").Value = 29999
$ws.Range("L51").Value = 29999
$ws.Range("N51").Value = -30967

$ws.Range("H70").Value = 1121.25
$ws.Range("I70").Value = 1162
$ws.Range("J70").Value = 999
$ws.Range("K70").Value = 3486
$ws.Range("L70").Value = 2997
$ws.Range("M70").Value = -3216
$ws.Range("N70").Value = -3537

$ws.Range("H73").Value = 1121.25
$ws.Range("I73").Value = 1162
$ws.Range("J73").Value = 999
$ws.Range("K73").Value = 3486
$ws.Range("L73").Value = 2997
$ws.Range("M73").Value = -2550
$ws.Range("N73").Value = -4869

$ws.Range("H98").Value = 582.6429000000001
$ws.Range("I98").Value = 614.0833
$ws.Range("J98").Value = 394
$ws.Range("K98").Value = 614.0833
$ws.Range("L98").Value = 394
$ws.Range("M98").Value = 883.9167
$ws.Range("N98").Value = -3390

$ws.Range("H122").Value = 582.6429000000001
$ws.Range("I122").Value = 614.0833
$ws.Range("J122").Value = 394
$ws.Range("K122").Value = 1842.2499
$ws.Range("L122").Value = 1182
$ws.Range("M122").Value = 607.7501
$ws.Range("N122").Value = -6082

$ws.Range("H131").Value = 5490.5713
$ws.Range("I131").Value = 5490.5713
$ws.Range("K131").Value = 16471.7139
$ws.Range("M131").Value = -11431.7139

$ws.Range("H137").Value = 2669.6667
$ws.Range("I137").Value = 1867.4546
$ws.Range("J137").Value = 6199.4
$ws.Range("K137").Value = 5602.3638
$ws.Range("L137").Value = 18598.2
$ws.Range("M137").Value = -3052.3638
$ws.Range("N137").Value = -23698.2

$ws.Range("H138").Value = 3368
$ws.Range("I138").Value = 1381.55
$ws.Range("K138").Value = 4144.65
$ws.Range("M138").Value = 995.3500000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2300.5833
$ws.Range("I61").Value = 2320.889
$ws.Range("K61").Value = 2320.889
$ws.Range("M61").Value = -2108.889

$ws.Range("H122").Value = 1166.4117
$ws.Range("I122").Value = 1166.4117
$ws.Range("K122").Value = 3499.2351
$ws.Range("M122").Value = -1049.2351

$ws.Range("H132").Value = 2045
$ws.Range("I132").Value = 1707.697
$ws.Range("K132").Value = 5123.090999999999
$ws.Range("M132").Value = -2593.090999999999

$ws.Range("H136").Value = 2300.5833
$ws.Range("I136").Value = 2320.889
$ws.Range("K136").Value = 6962.667
$ws.Range("M136").Value = -4412.667

$ws.Range("H138").Value = 90000
$ws.Range("J138").Value = 90000
$ws.Range("L138").Value = 90000
$ws.Range("N138").Value = -100280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2258.647
$ws.Range("I134").Value = 2231.6428
$ws.Range("K134").Value = 6694.928400000001
$ws.Range("M134").Value = -4159.928400000001

$ws.Range("H135").Value = 46153.1
$ws.Range("J135").Value = 46153.1
$ws.Range("L135").Value = 46153.1
$ws.Range("N135").Value = -56293.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4348.9287
$ws.Range("I31").Value = 2599
$ws.Range("J31").Value = 7498.8
$ws.Range("K31").Value = 2599
$ws.Range("L31").Value = 7498.8
$ws.Range("M31").Value = -2304
$ws.Range("N31").Value = -8088.8

$ws.Range("H34").Value = 4348.9287
$ws.Range("I34").Value = 2599
$ws.Range("J34").Value = 7498.8
$ws.Range("K34").Value = 2599
$ws.Range("L34").Value = 7498.8
$ws.Range("M34").Value = -2397
$ws.Range("N34").Value = -7902.8

$ws.Range("H58").Value = 2494
$ws.Range("I58").Value = 2208.8
$ws.Range("K58").Value = 2208.8
$ws.Range("M58").Value = -2005.8

$ws.Range("H107").Value = 907.0833
$ws.Range("I107").Value = 611
$ws.Range("J107").Value = 1499.25
$ws.Range("K107").Value = 611
$ws.Range("L107").Value = 1499.25
$ws.Range("M107").Value = 1309
$ws.Range("N107").Value = -5339.25

$ws.Range("H132").Value = 2924.4783
$ws.Range("I132").Value = 2830.775
$ws.Range("K132").Value = 8492.325000000001
$ws.Range("M132").Value = -5962.325000000001

$ws.Range("H134").Value = 4242.375
$ws.Range("I134").Value = 4643.4614
$ws.Range("J134").Value = 2504.3333
$ws.Range("K134").Value = 13930.3842
$ws.Range("L134").Value = 7512.999899999999
$ws.Range("M134").Value = -11395.3842
$ws.Range("N134").Value = -12582.9999

$ws.Range("H136").Value = 2494
$ws.Range("I136").Value = 2208.8
$ws.Range("K136").Value = 6626.400000000001
$ws.Range("M136").Value = -4076.400000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1832.1666
$ws.Range("J107").Value = 198.6
$ws.Range("L107").Value = 595.8
$ws.Range("N107").Value = -4435.8

$ws.Range("H134").Value = 2075.3333
$ws.Range("I134").Value = 2075.3333
$ws.Range("K134").Value = 6225.999899999999
$ws.Range("M134").Value = -1155.999899999999

$ws.Range("H139").Value = 3481.5
$ws.Range("J139").Value = 5033
$ws.Range("L139").Value = 15099
$ws.Range("N139").Value = -25379

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4933.3335
$ws.Range("I102").Value = 4933.3335
$ws.Range("K102").Value = 4933.3335
$ws.Range("M102").Value = -3311.3335

$ws.Range("H132").Value = 3002.4666
$ws.Range("I132").Value = 2670
$ws.Range("K132").Value = 8010
$ws.Range("M132").Value = -5480

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2025.5
$ws.Range("I93").Value = 2025.5
$ws.Range("K93").Value = 2025.5
$ws.Range("M93").Value = -777.5

$ws.Range("H136").Value = 5202.8
$ws.Range("I136").Value = 5202.8
$ws.Range("K136").Value = 15608.4
$ws.Range("M136").Value = -13058.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1187.5
$ws.Range("J100").Value = 75
$ws.Range("L100").Value = 150
$ws.Range("N100").Value = -1232

$ws.Range("H107").Value = 997.5
$ws.Range("I107").Value = 997.6667
$ws.Range("J107").Value = 997
$ws.Range("K107").Value = 2993.0001
$ws.Range("L107").Value = 2991
$ws.Range("M107").Value = -1073.0001
$ws.Range("N107").Value = -6831

$ws.Range("H113").Value = 532.5714
$ws.Range("I113").Value = 532.5714
$ws.Range("K113").Value = 1597.7142
$ws.Range("M113").Value = 572.2857999999999

$ws.Range("H136").Value = 1015.7
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
